$d = $word.ActiveDocument

# Update the date heading
[void]$d.Content.Find.Execute("2025-02-02 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-02-03 Monday", 2)

# Update the division problems in the table, cell by cell (row, column)
# so that duplicate problem text (e.g. "12÷6=") is handled unambiguously.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "53÷8="
$t.Cell(1,2).Range.Text  = "40÷8="
$t.Cell(1,3).Range.Text  = "72÷2="
$t.Cell(1,4).Range.Text  = "35÷8="
$t.Cell(1,5).Range.Text  = "11÷9="

$t.Cell(5,1).Range.Text  = "53÷4="
$t.Cell(5,2).Range.Text  = "95÷5="
$t.Cell(5,3).Range.Text  = "68÷8="
$t.Cell(5,4).Range.Text  = "51÷2="
$t.Cell(5,5).Range.Text  = "34÷4="

$t.Cell(9,1).Range.Text  = "96÷5="
$t.Cell(9,2).Range.Text  = "49÷5="
$t.Cell(9,3).Range.Text  = "74÷7="
$t.Cell(9,4).Range.Text  = "95÷8="
$t.Cell(9,5).Range.Text  = "11÷4="

$t.Cell(13,1).Range.Text = "46÷4="
$t.Cell(13,2).Range.Text = "66÷6="
$t.Cell(13,3).Range.Text = "40÷3="
$t.Cell(13,4).Range.Text = "27÷2="
$t.Cell(13,5).Range.Text = "87÷9="

$t.Cell(17,1).Range.Text = "25÷6="
$t.Cell(17,2).Range.Text = "98÷9="
$t.Cell(17,3).Range.Text = "67÷8="
$t.Cell(17,4).Range.Text = "99÷8="
$t.Cell(17,5).Range.Text = "41÷9="
